# Clean horsepower text strings in column D ("Power Output") into plain
# numeric values. Rows that previously held "X hp" / "X,XXX hp" style text
# now hold the bare number so the sheet can do numeric work with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$powerValues = @{
    2  = 150
    3  = 150
    4  = 150
    5  = 150
    6  = 400
    7  = 300
    18 = 580
    20 = 800
    21 = 940
    22 = 600
    25 = 600
    26 = 600
    27 = 1000
    28 = 1100
}

foreach ($row in $powerValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $powerValues[$row]
}

# The 4-digit horsepower figures (1,000 / 1,100) keep a thousands-separator
# number format so they still read the same way they did as text.
$ws.Range("D27:D28").NumberFormat = "#,##0"

# Reflect the last active selection from the authored workbook.
$ws.Range("J19").Select()
